# FreeCrmTestData.xlsx - refresh the "contacts" sheet test data.
#
# Before:
#   title firstname lastname company
#   Mr.   Tom       Peter    Google
#   Dr.   David     Cris     Amazon
#   Mrs.  Mukta     Sharma   Ebay
#
# After:
#   firstname lastname company
#   Marry     Ds       Ebay
#   David     Thomas   Home Depot

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("contacts")

# Drop the "title" column entirely - firstname/lastname/company shift left.
$ws.Columns.Item(1).Delete() | Out-Null

# Drop the third data row (old row 4: Mrs./Mukta/Sharma/Ebay).
$ws.Rows.Item(4).Delete() | Out-Null

# Replace the remaining two data rows with the new sample values.
$ws.Range("A2").Value = "Marry"
$ws.Range("B2").Value = "Ds"
$ws.Range("C2").Value = "Ebay"

$ws.Range("A3").Value = "David"
$ws.Range("B3").Value = "Thomas"
$ws.Range("C3").Value = "Home Depot"

# Leave the cursor where the saved file shows it.
$ws.Range("B5").Select() | Out-Null
